# household.xlsx regen: update settings!form_version and the saved
# selection/active-sheet UI state to match the author's last interaction.

$wb = $excel.ActiveWorkbook

# --- settings sheet: form_version 1 -> 20130408 -----------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Value = 20130408

# --- UI state: the workbook was last saved with "settings" active and -
# cell B7 selected there (previously "choices" was active/tabSelected). -
$settings.Activate()
$settings.Range("B7").Select()
